$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update company name for row 3
$ws.Range("B3").Value = "Kino Co., Ltd. (GTSM:4154)"

# Values that are identical between row 2 and row 3
$colValues = @{
    "D"  = 0.0542
    "G"  = -0.05412844036697247
    "H"  = -0.05458715596330275
    "I"  = -0.03738532110091743
    "J"  = -0.03738532110091743
    "K"  = -0.726
    "L"  = -0.03330275229357798
    "U"  = 6.75
    "V"  = 0.2934782608695652
    "W"  = -0.09797570850202429
    "X"  = 0.05003445640734468
    "Y"  = -0.148010164909369
    "Z"  = 3.679324894514768
    "AA" = -0.1375527426160338
    "AB" = 0.04615367772445398
    "AC" = -0.1837064203404877
    "AD" = 4.13
    "AE" = 0
    "AF" = 4.13
    "AG" = -2.62
    "AH" = 0.1522300036859565
    "AI" = 0.3005822416302766
    "AJ" = -0.1285574092247301
    "AK" = -0.374821173104435
    "AL" = 0.137
    "AM" = 0.114
    "AN" = -6.154992548435171
    "AO" = -5.948905109489051
    "AP" = 3.904619970193741
    "AQ" = -7.149122807017542
}

foreach ($col in $colValues.Keys) {
    $value = $colValues[$col]
    $ws.Range("$col" + "2").Value = $value
    $ws.Range("$col" + "3").Value = $value
}
